$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 89999
$ws.Range("J87").Value = 89999
$ws.Range("L87").Value = 89999
$ws.Range("N87").Value = -92495
$ws.Range("H90").Value = 89999
$ws.Range("J90").Value = 89999
$ws.Range("L90").Value = 269997
$ws.Range("N90").Value = -282477
$ws.Range("H97").Value = 4384
$ws.Range("J97").Value = 3375
$ws.Range("L97").Value = 10125
$ws.Range("N97").Value = -11117
$ws.Range("H100").Value = 605.55
$ws.Range("I100").Value = 538.8125
$ws.Range("K100").Value = 538.8125
$ws.Range("M100").Value = 2.1875
$ws.Range("H125").Value = 817.5
$ws.Range("I125").Value = 806.25
$ws.Range("J125").Value = 862.5
$ws.Range("K125").Value = 7256.25
$ws.Range("L125").Value = 7762.5
$ws.Range("M125").Value = -4796.25
$ws.Range("N125").Value = -12682.5
$ws.Range("H138").Value = 14843.223
$ws.Range("I138").Value = 7948.75
$ws.Range("K138").Value = 23846.25
$ws.Range("M138").Value = -18706.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13221.419
$ws.Range("I32").Value = 5368
$ws.Range("J32").Value = 26892.186
$ws.Range("K32").Value = 5368
$ws.Range("L32").Value = 26892.186
$ws.Range("M32").Value = -5081
$ws.Range("N32").Value = -27466.186
$ws.Range("H61").Value = 1341.8572
$ws.Range("I61").Value = 1354.2693
$ws.Range("K61").Value = 1354.2693
$ws.Range("M61").Value = -1142.2693
$ws.Range("H119").Value = 25000
$ws.Range("J119").Value = 25000
$ws.Range("L119").Value = 25000
$ws.Range("N119").Value = -34676
$ws.Range("H122").Value = 419591.4
$ws.Range("J122").Value = 4609.5
$ws.Range("L122").Value = 13828.5
$ws.Range("N122").Value = -18728.5
$ws.Range("H136").Value = 1341.8572
$ws.Range("I136").Value = 1354.2693
$ws.Range("K136").Value = 4062.8079
$ws.Range("M136").Value = -1512.8079

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3368.2856
$ws.Range("I20").Value = 2992.6
$ws.Range("J20").Value = 4307.5
$ws.Range("K20").Value = 2992.6
$ws.Range("L20").Value = 4307.5
$ws.Range("M20").Value = -2745.6
$ws.Range("N20").Value = -4801.5
$ws.Range("H64").Value = 1217.091
$ws.Range("I64").Value = 1068.6
$ws.Range("J64").Value = 1340.8334
$ws.Range("K64").Value = 1068.6
$ws.Range("L64").Value = 1340.8334
$ws.Range("M64").Value = -843.5999999999999
$ws.Range("N64").Value = -1790.8334
$ws.Range("H67").Value = 1217.091
$ws.Range("I67").Value = 1068.6
$ws.Range("J67").Value = 1340.8334
$ws.Range("K67").Value = 1068.6
$ws.Range("L67").Value = 1340.8334
$ws.Range("M67").Value = -288.5999999999999
$ws.Range("N67").Value = -2900.8334
$ws.Range("H105").Value = 4258.6772
$ws.Range("I105").Value = 3612.2222
$ws.Range("K105").Value = 3612.2222
$ws.Range("M105").Value = -1865.2222

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3610.652
$ws.Range("I31").Value = 1171.4546
$ws.Range("K31").Value = 1171.4546
$ws.Range("M31").Value = -876.4546
$ws.Range("H34").Value = 3610.652
$ws.Range("I34").Value = 1171.4546
$ws.Range("K34").Value = 1171.4546
$ws.Range("M34").Value = -969.4546
$ws.Range("H134").Value = 4781.9287
$ws.Range("I134").Value = 3908.818
$ws.Range("K134").Value = 11726.454
$ws.Range("M134").Value = -9191.454000000002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 22.846153
$ws.Range("I12").Value = 25.857143
$ws.Range("K12").Value = 77.57142899999999
$ws.Range("M12").Value = 95.42857100000001
$ws.Range("H81").Value = 2018.8334
$ws.Range("I81").Value = 600
$ws.Range("J81").Value = 2728.25
$ws.Range("K81").Value = 1800
$ws.Range("L81").Value = 8184.75
$ws.Range("M81").Value = -677
$ws.Range("N81").Value = -10430.75
$ws.Range("H84").Value = 2018.8334
$ws.Range("I84").Value = 600
$ws.Range("J84").Value = 2728.25
$ws.Range("K84").Value = 5400
$ws.Range("L84").Value = 24554.25
$ws.Range("M84").Value = 216
$ws.Range("N84").Value = -35786.25
$ws.Range("H120").Value = 14377.223
$ws.Range("H121").Value = 2047.625
$ws.Range("J121").Value = 3673.25
$ws.Range("L121").Value = 11019.75
$ws.Range("N121").Value = -13639.75
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 3000
$ws.Range("N125").ClearContents()
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = 1920

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 522.5714
$ws.Range("J2").Value = 882.2857
$ws.Range("L2").Value = 882.2857
$ws.Range("N2").Value = -1108.2857
$ws.Range("H70").Value = 7332.222
$ws.Range("J70").Value = 7332.222
$ws.Range("L70").Value = 7332.222
$ws.Range("N70").Value = -7872.222
$ws.Range("H73").Value = 7332.222
$ws.Range("J73").Value = 7332.222
$ws.Range("L73").Value = 7332.222
$ws.Range("N73").Value = -9204.222
$ws.Range("H97").Value = 2149.44
$ws.Range("I97").Value = 2214.6316
$ws.Range("J97").Value = 1943
$ws.Range("K97").Value = 2214.6316
$ws.Range("L97").Value = 1943
$ws.Range("M97").Value = -1718.6316
$ws.Range("N97").Value = -2935
$ws.Range("H102").Value = 2244
$ws.Range("I102").Value = 1062.5
$ws.Range("J102").Value = 3294.2222
$ws.Range("K102").Value = 1062.5
$ws.Range("L102").Value = 3294.2222
$ws.Range("M102").Value = 559.5
$ws.Range("N102").Value = -6538.2222
$ws.Range("H122").Value = 357880.53
$ws.Range("I122").Value = 93094.73
$ws.Range("J122").Value = 503512.7
$ws.Range("K122").Value = 279284.19
$ws.Range("L122").Value = 1510538.1
$ws.Range("M122").Value = -276834.19
$ws.Range("N122").Value = -1515438.1

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 6700099.5
$ws.Range("J2").Value = 50149.5
$ws.Range("L2").Value = 50149.5
$ws.Range("N2").Value = -50373.5
$ws.Range("H22").Value = 1996.8572
$ws.Range("I22").Value = 1996.6666
$ws.Range("K22").Value = 1996.6666
$ws.Range("M22").Value = -1701.6666
$ws.Range("H27").Value = 1996.8572
$ws.Range("I27").Value = 1996.6666
$ws.Range("K27").Value = 1996.6666
$ws.Range("M27").Value = -1889.6666
$ws.Range("H46").Value = 2499.1155
$ws.Range("I46").Value = 1577.7368
$ws.Range("K46").Value = 1577.7368
$ws.Range("M46").Value = -1389.7368
$ws.Range("H82").Value = 1969.8235
$ws.Range("I82").Value = 2061.0908
$ws.Range("J82").Value = 1802.5
$ws.Range("K82").Value = 2061.0908
$ws.Range("L82").Value = 1802.5
$ws.Range("M82").Value = -1700.0908
$ws.Range("N82").Value = -2524.5
$ws.Range("H85").Value = 1969.8235
$ws.Range("I85").Value = 2061.0908
$ws.Range("J85").Value = 1802.5
$ws.Range("K85").Value = 2061.0908
$ws.Range("L85").Value = 1802.5
$ws.Range("M85").Value = -813.0907999999999
$ws.Range("N85").Value = -4298.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6832.05
$ws.Range("I62").Value = 3632.3333
$ws.Range("J62").Value = 7396.706
$ws.Range("K62").Value = 3632.3333
$ws.Range("L62").Value = 7396.706
$ws.Range("M62").Value = -3008.3333
$ws.Range("N62").Value = -8644.706
$ws.Range("H65").Value = 6832.05
$ws.Range("I65").Value = 3632.3333
$ws.Range("J65").Value = 7396.706
$ws.Range("K65").Value = 18161.6665
$ws.Range("L65").Value = 36983.53
$ws.Range("M65").Value = -15041.6665
$ws.Range("N65").Value = -43223.53
$ws.Range("H119").Value = 40749.75
$ws.Range("J119").Value = 40749.75
$ws.Range("L119").Value = 40749.75
$ws.Range("N119").Value = -50425.75
$ws.Range("H122").Value = 2007.9
$ws.Range("J122").Value = 1850
$ws.Range("L122").Value = 5550
$ws.Range("N122").Value = -10450
